$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "AutoOpp"
$ws.Range("A3").Value = "AutoOpp"

$ws.Range("A4").Select()
